$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (D) and 1h volume change (E) columns
# Values that look numeric need the cell pre-formatted as Text so Excel
# keeps the exact display string (matching dotted price formatting, e.g. "596.89")
# instead of silently parsing them into a float.
$textCells = @("D5", "D6", "D8", "D10", "D14", "D16", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D27", "D31", "D32", "D36", "D38", "D41", "D42", "D43", "D46", "D47", "D48", "D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '68.381.90'
$ws.Range('E2').Value = '  -0.04%  '
$ws.Range('D3').Value = '2.646.07'
$ws.Range('E3').Value = '  -0.03%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = '596.89'
$ws.Range('E5').Value = '  -0.27%  '
$ws.Range('D6').Value = '158.86'
$ws.Range('E6').Value = '  +2.65%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').Value = '0.538'
$ws.Range('E8').Value = '  -1.60%  '
$ws.Range('D9').Value = '2.645.37'
$ws.Range('D10').Value = '0.143'
$ws.Range('E10').Value = '  -1.67%  '
$ws.Range('E12').Value = '  +0.14%  '
$ws.Range('D14').Value = '27.91'
$ws.Range('E14').Value = '  -0.61%  '
$ws.Range('D15').Value = '3.131.20'
$ws.Range('E15').Value = '  +0.08%  '
$ws.Range('D16').Value = '0.0000187'
$ws.Range('E16').Value = '  -3.09%  '
$ws.Range('D17').Value = '68.643.55'
$ws.Range('E17').Value = '  +0.48%  '
$ws.Range('D18').Value = '2.608.59'
$ws.Range('E18').Value = '  -1.81%  '
$ws.Range('D19').Value = '11.41'
$ws.Range('E19').Value = '  +0.43%  '
$ws.Range('D20').Value = '362.63'
$ws.Range('E20').Value = '  -0.40%  '
$ws.Range('D21').Value = '7.38'
$ws.Range('E21').Value = '  -1.58%  '
$ws.Range('D22').Value = '4.40'
$ws.Range('E22').Value = '  +0.48%  '
$ws.Range('D23').Value = '4.77'
$ws.Range('E23').Value = '  -2.53%  '
$ws.Range('D24').Value = '2.07'
$ws.Range('E24').Value = '  +0.16%  '
$ws.Range('D25').Value = '74.42'
$ws.Range('E25').Value = '  -0.55%  '
$ws.Range('E26').Value = '  +0.03%  '
$ws.Range('D27').Value = '9.78'
$ws.Range('E27').Value = '  -0.53%  '
$ws.Range('D28').Value = '2.778.57'
$ws.Range('E28').Value = '  +0.08%  '
$ws.Range('E29').Value = '  -3.29%  '
$ws.Range('E30').Value = '  -0.13%  '
$ws.Range('D31').Value = '562.80'
$ws.Range('E31').Value = '  -1.99%  '
$ws.Range('D32').Value = '8.05'
$ws.Range('E32').Value = '  -0.37%  '
$ws.Range('E33').Value = '  -1.99%  '
$ws.Range('E34').Value = '  -0.92%  '
$ws.Range('E35').Value = '  +3.63%  '
$ws.Range('D36').Value = '0.999'
$ws.Range('E36').Value = '  +0.01%  '
$ws.Range('E37').Value = '  -1.76%  '
$ws.Range('D38').Value = '159.59'
$ws.Range('E38').Value = '  -0.73%  '
$ws.Range('E39').Value = '  +1.22%  '
$ws.Range('E40').Value = '  -1.19%  '
$ws.Range('D41').Value = '1.86'
$ws.Range('E41').Value = '  -1.59%  '
$ws.Range('D42').Value = '5.32'
$ws.Range('E42').Value = '  -1.18%  '
$ws.Range('D43').Value = '2.62'
$ws.Range('E43').Value = '  -1.76%  '
$ws.Range('D44').Value = '0.0₆0319'
$ws.Range('E44').Value = '  -5.33%  '
$ws.Range('D46').Value = '157.76'
$ws.Range('E46').Value = '  +0.59%  '
$ws.Range('D47').Value = '3.81'
$ws.Range('E47').Value = '  +1.33%  '
$ws.Range('D48').Value = '21.97'
$ws.Range('E48').Value = '  +0.30%  '
$ws.Range('E49').Value = '  -1.52%  '
$ws.Range('D50').Value = '0.0774'
$ws.Range('E50').Value = '  -1.70%  '
$ws.Range('E51').Value = '  +0.98%  '
